# Adds 2015 NRI data as a new column I, and refreshes the underlying
# weighted-mean values for columns B:H on rows 8-20 (lccL7_pcnt through
# lccL78_pcnt) to reflect the recomputed series.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Give the new column I the same number format / border style as
#     column H before populating any values, by copying its formatting. ---
$ws.Range("H1:H26").Copy()
$ws.Range("I1:I26").PasteSpecial(-4122)

# --- Updated B:H values for rows 8-20 ---
$updatedRows = @{
    8 = @(16.768915176391602,16.774328231811523,16.765037536621094,16.766012191772461,16.733757019042969,16.722023010253906,16.715496063232422)
    9 = @(1.8552695512771606,1.8545185327529907,1.8525469303131104,1.8450156450271606,1.8361639976501465,1.8287836313247681,1.8227254152297974)
    10 = @(18.196645736694336,18.185880661010742,18.167137145996094,18.158349990844727,18.149513244628906,18.124229431152344,18.123319625854492)
    11 = @(17.783185958862305,17.799587249755859,17.817337036132813,17.83479118347168,17.855457305908203,17.877235412597656,17.882226943969727)
    12 = @(12.014049530029297,12.01760196685791,12.023721694946289,12.017516136169434,12.02822208404541,12.029709815979004,12.032098770141602)
    13 = @(2.0847756862640381,2.0832056999206543,2.0810856819152832,2.0796327590942383,2.0810973644256592,2.0805275440216064,2.081218957901001)
    14 = @(15.732243537902832,15.722705841064453,15.72838020324707,15.723287582397461,15.737751960754395,15.755236625671387,15.758761405944824)
    15 = @(14.565268516540527,14.5625,14.5667724609375,14.577980041503906,14.583365440368652,14.588701248168945,14.59078311920166)
    16 = @(0.99964636564254761,0.99967110157012939,0.99798119068145752,0.99741446971893311,0.99466949701309204,0.99355357885360718,0.99336904287338257)
    17 = @(20.051916122436523,20.040399551391602,20.019683837890625,20.003364562988281,19.985677719116211,19.953012466430664,19.946044921875)
    18 = @(29.797235488891602,29.817188262939453,29.841058731079102,29.85230827331543,29.88368034362793,29.906944274902344,29.914325714111328)
    19 = @(17.817018508911133,17.805912017822266,17.809465408325195,17.802921295166016,17.818849563598633,17.835763931274414,17.839981079101563)
    20 = @(15.564914703369141,15.56217098236084,15.564753532409668,15.575394630432129,15.578034400939941,15.582255363464355,15.584152221679688)
}

foreach ($r in $updatedRows.Keys) {
    $rowVals = $updatedRows[$r]
    for ($i = 0; $i -lt 7; $i++) {
        $ws.Cells.Item($r, $i + 2).Value = $rowVals[$i]
    }
}

# --- New column I (2015) values. Rows 21, 22 and 23 stay blank (only
#     formatting carried over above), matching the source data availability. ---
$newCol = @{
    1 = 2015
    2 = 25.903202056884766
    3 = 1.2667135000228882
    4 = 29.242368698120117
    5 = 8.5734872817993164
    6 = 28.506320953369141
    7 = 6.5079078674316406
    8 = 16.715444564819336
    9 = 1.8197025060653687
    10 = 18.123176574707031
    11 = 17.884641647338867
    12 = 12.033391952514648
    13 = 2.0811939239501953
    14 = 15.755876541137695
    15 = 14.59307861328125
    16 = 0.99349385499954224
    17 = 19.942878723144531
    18 = 29.918033599853516
    19 = 17.837070465087891
    20 = 15.586572647094727
    24 = 61.171245574951172
    25 = 16.908662796020508
    26 = 16.908662796020508
}

foreach ($r in $newCol.Keys) {
    $ws.Cells.Item($r, 9).Value = $newCol[$r]
}
